# Add "Promotions" and "Groups" columns to the exam parameters table,
# located right after the "oral/written" column and before "start date".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lo = $ws.ListObjects.Item(1)
$tableName = $lo.Name
$tableStyle = "TableStyleLight8"

# 1. Insert two blank columns at F:G. This shifts the existing
#    "start date"/"final date"/"is weekend ok?" (and the helper yes/no & oral/written
#    list values) columns from F:I to H:K, and Excel automatically keeps the
#    data validation target ranges (sqref) in sync.
$ws.Range("F1:G1").EntireColumn.Insert()

# 2. The data-validation list formulas still point at the old column (I); repoint
#    them at the new helper column (K).
$ws.Range("E2:E7").Validation.Formula1 = "=`$K`$4:`$K`$5"
$ws.Range("J2:J7").Validation.Formula1 = "=`$K`$2:`$K`$3"

# 3. Fill in the new "Promotions" / "Groups" columns.
$ws.Range("F1").Value = "Promotions"
$ws.Range("G1").Value = "Groups"

for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 6).Value = 172
}
$ws.Range("G2").Value = "BHK=pilot"
$ws.Range("G4").Value = "BHK=navy"

# 3b. Match the column widths used for the surrounding columns (E inherited its
#     width onto the two new columns, and "start date" ends up a bit narrower).
$ws.Range("F1:G1").ColumnWidth = 14.5
$ws.Range("H1").ColumnWidth = 10.333333333333334

# 4. Convert the existing table to a normal range (keeps the data/formatting,
#    unlike Delete which would clear the cells) so a new table definition can
#    be rebuilt with the correct column layout.
$lo.Unlist()

# 5. Re-assert every header cell value so the new table picks up proper
#    column names for all columns (not just the two we just touched).
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Unavailability"
$ws.Range("C1").Value = "Amount days"
$ws.Range("D1").Value = "preparation days"
$ws.Range("E1").Value = "oral/written"
$ws.Range("F1").Value = "Promotions"
$ws.Range("G1").Value = "Groups"
$ws.Range("H1").Value = "start date"
$ws.Range("I1").Value = "final date"
$ws.Range("J1").Value = "is weekend ok?"

# 6. Rebuild the table over the expanded range and restore its name/style.
$lo2 = $ws.ListObjects.Add(1, $ws.Range("A1:J7"), [System.Reflection.Missing]::Value, 1)
$lo2.Name = $tableName
$lo2.TableStyle = $tableStyle

# 7. Leave the selection on the newly filled-in "BHK=navy" cell, as in the authored edit.
[void]$ws.Range("G4").Select()
